$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# Header text updates (report volume/number and date range), new weekly
# crime data collected for the week of 11/6/2023 - 11/12/2023.
# --------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --------------------------------------------------------------------
# Cells that change data type from text ("0") to a real number.
# Copy the number format from an untouched numeric cell in the same
# row/style family so the resulting style matches a genuine numeric
# cell (not just a coerced text cell).
# --------------------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = $ws.Range("D15").NumberFormat

$ws.Range("C23").Value = 2
$ws.Range("C23").NumberFormat = $ws.Range("G23").NumberFormat

$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = $ws.Range("G30").NumberFormat

# --------------------------------------------------------------------
# Cells that change data type from a number to text ("0" or "***.*").
# Use Copy/PasteSpecial from a cell that already holds that exact text
# (and style) elsewhere on the sheet so both the shared-string value
# and the cell style end up correct.
# --------------------------------------------------------------------
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

$ws.Range("D14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null

$ws.Range("E14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null

# --------------------------------------------------------------------
# Remaining plain numeric value updates (same type/style, new figures).
# --------------------------------------------------------------------
$ws.Range("N14").Value = -9.523809523809
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 32
$ws.Range("J15").Value = 32
$ws.Range("L15").Value = -21.951219512195
$ws.Range("M15").Value = -37.254901960784
$ws.Range("N15").Value = -54.929577464788
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 270
$ws.Range("J16").Value = 234
$ws.Range("K16").Value = 15.384615384615
$ws.Range("L16").Value = 51.685393258427
$ws.Range("M16").Value = -27.027027027027
$ws.Range("N16").Value = -77.649006622516
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -6.25
$ws.Range("F17").Value = 73
$ws.Range("G17").Value = 68
$ws.Range("H17").Value = 7.352941176470
$ws.Range("I17").Value = 849
$ws.Range("J17").Value = 697
$ws.Range("K17").Value = 21.807747489239
$ws.Range("L17").Value = 62.643678160919
$ws.Range("M17").Value = 96.983758700696
$ws.Range("N17").Value = -20.356472795497
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -36.363636363636
$ws.Range("F18").Value = 30
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 277
$ws.Range("K18").Value = 8.303249097472
$ws.Range("L18").Value = 41.509433962264
$ws.Range("M18").Value = -45.355191256830
$ws.Range("N18").Value = -89.973262032085
$ws.Range("C19").Value = 42
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 136
$ws.Range("G19").Value = 134
$ws.Range("H19").Value = 1.492537313432
$ws.Range("I19").Value = 1407
$ws.Range("J19").Value = 1214
$ws.Range("K19").Value = 15.897858319604
$ws.Range("L19").Value = 49.680851063829
$ws.Range("M19").Value = 67.102137767220
$ws.Range("N19").Value = -1.951219512195
$ws.Range("C20").Value = 7
$ws.Range("E20").Value = -65
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 66
$ws.Range("H20").Value = -46.969696969697
$ws.Range("I20").Value = 403
$ws.Range("J20").Value = 431
$ws.Range("K20").Value = -6.496519721577
$ws.Range("L20").Value = 78.318584070796
$ws.Range("M20").Value = 30.42071197411
$ws.Range("N20").Value = -90.484061393152
$ws.Range("C21").Value = 79
$ws.Range("D21").Value = 80
$ws.Range("E21").Value = -1.25
$ws.Range("F21").Value = 298
$ws.Range("G21").Value = 320
$ws.Range("H21").Value = -6.875
$ws.Range("I21").Value = 3280
$ws.Range("J21").Value = 2894
$ws.Range("K21").Value = 13.337940566689
$ws.Range("L21").Value = 53.918348193336
$ws.Range("M21").Value = 27.825409197194
$ws.Range("N21").Value = -70.257526296699
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 20
$ws.Range("I23").Value = 118
$ws.Range("J23").Value = 84
$ws.Range("K23").Value = 40.476190476190
$ws.Range("L23").Value = 25.531914893617
$ws.Range("M23").Value = 87.301587301587
$ws.Range("C24").Value = 92
$ws.Range("D24").Value = 91
$ws.Range("E24").Value = 1.098901098901
$ws.Range("G24").Value = 339
$ws.Range("H24").Value = 13.274336283185
$ws.Range("I24").Value = 3759
$ws.Range("J24").Value = 3419
$ws.Range("K24").Value = 9.944428195378
$ws.Range("L24").Value = 73.305670816044
$ws.Range("M24").Value = 14.324817518248
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 125
$ws.Range("G25").Value = 117
$ws.Range("H25").Value = 6.837606837606
$ws.Range("I25").Value = 1528
$ws.Range("J25").Value = 1433
$ws.Range("K25").Value = 6.629448709002
$ws.Range("L25").Value = 30.264279624893
$ws.Range("M25").Value = -14.732142857142
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 59
$ws.Range("J26").Value = 66
$ws.Range("K26").Value = -10.606060606060
$ws.Range("L26").Value = -25.316455696202
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 20
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 76.923076923076
$ws.Range("I27").Value = 189
$ws.Range("J27").Value = 153
$ws.Range("K27").Value = 23.529411764705
$ws.Range("L27").Value = 44.274809160305
$ws.Range("N28").Value = -72.549019607843
$ws.Range("N29").Value = -71.111111111111
$ws.Range("F30").Value = 3
$ws.Range("H30").Value = 200
$ws.Range("I30").Value = 18
$ws.Range("K30").Value = -18.181818181818
$ws.Range("L30").Value = 5.882352941176
